$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 31 (ALC)
$ws.Range("H31").Value = 1490.375
$ws.Range("I31").Value = 132
$ws.Range("J31").Value = 10999
$ws.Range("K31").Value = 396
$ws.Range("L31").Value = 32997
$ws.Range("M31").Value = -166
$ws.Range("N31").Value = -33457

# Row 98 (ALC)
$ws.Range("H98").Value = 3481.647
$ws.Range("I98").Value = 1476.6154
$ws.Range("J98").Value = 9998
$ws.Range("K98").Value = 1476.6154
$ws.Range("L98").Value = 9998
$ws.Range("M98").Value = 21.38460000000009
$ws.Range("N98").Value = -12994

# Row 116 (ALC)
$ws.Range("H116").Value = 5284.864
$ws.Range("I116").Value = 5068.4
$ws.Range("K116").Value = 5068.4
$ws.Range("M116").Value = -1626.4

# Row 122 (ALC)
$ws.Range("H122").Value = 3481.647
$ws.Range("I122").Value = 1476.6154
$ws.Range("J122").Value = 9998
$ws.Range("K122").Value = 4429.8462
$ws.Range("L122").Value = 29994
$ws.Range("M122").Value = -1979.8462
$ws.Range("N122").Value = -34894

# Row 137 (ALC)
$ws.Range("H137").Value = 2435.7273
$ws.Range("I137").Value = 1643.6666
$ws.Range("J137").Value = 6000
$ws.Range("K137").Value = 4930.9998
$ws.Range("L137").Value = 18000
$ws.Range("M137").Value = -2380.9998
$ws.Range("N137").Value = -23100

# Row 138 (ALC)
$ws.Range("H138").Value = 3141.4443
$ws.Range("J138").Value = 5451.7144
$ws.Range("L138").Value = 16355.1432
$ws.Range("N138").Value = -26635.1432

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (ARM)
$ws.Range("H2").Value = 736572.6
$ws.Range("I2").Value = 1051531.8
$ws.Range("J2").Value = 1668.1666
$ws.Range("K2").Value = 1051531.8
$ws.Range("L2").Value = 1668.1666
$ws.Range("M2").Value = -1051418.8
$ws.Range("N2").Value = -1894.1666

# Row 10 (ARM)
$ws.Range("H10").Value = 259500
$ws.Range("I10").Value = 12666.667
$ws.Range("J10").Value = 1000000
$ws.Range("K10").Value = 12666.667
$ws.Range("L10").Value = 1000000
$ws.Range("M10").Value = -12496.667
$ws.Range("N10").Value = -1000340

# Row 61 (ARM)
$ws.Range("H61").Value = 50001852
$ws.Range("I61").Value = 55557390
$ws.Range("K61").Value = 55557390
$ws.Range("M61").Value = -55557178

# Row 74 (ARM)
$ws.Range("H74").Value = 47625460
$ws.Range("I74").Value = 47625460
$ws.Range("K74").Value = 47625460
$ws.Range("M74").Value = -47624586

# Row 77 (ARM)
$ws.Range("H77").Value = 47625460
$ws.Range("I77").Value = 47625460
$ws.Range("K77").Value = 238127300
$ws.Range("M77").Value = -238122932

# Row 110 (ARM)
$ws.Range("H110").Value = 144170.72
$ws.Range("I110").Value = 167882.5
$ws.Range("J110").Value = 1900
$ws.Range("K110").Value = 167882.5
$ws.Range("L110").Value = 1900
$ws.Range("M110").Value = -165837.5
$ws.Range("N110").Value = -5990

# Row 116 (ARM)
$ws.Range("H116").Value = 736572.6
$ws.Range("I116").Value = 1051531.8
$ws.Range("J116").Value = 1668.1666
$ws.Range("K116").Value = 1051531.8
$ws.Range("L116").Value = 1668.1666
$ws.Range("M116").Value = -1049237.8
$ws.Range("N116").Value = -6256.1666

# Row 122 (ARM)
$ws.Range("H122").Value = 3880.5557
$ws.Range("I122").Value = 2808
$ws.Range("K122").Value = 8424
$ws.Range("M122").Value = -5974

# Row 136 (ARM)
$ws.Range("H136").Value = 50001852
$ws.Range("I136").Value = 55557390
$ws.Range("K136").Value = 166672170
$ws.Range("M136").Value = -166669620

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (BSM)
$ws.Range("H3").Value = 736572.6
$ws.Range("I3").Value = 1051531.8
$ws.Range("J3").Value = 1668.1666
$ws.Range("K3").Value = 1051531.8
$ws.Range("L3").Value = 1668.1666
$ws.Range("M3").Value = -1051417.8
$ws.Range("N3").Value = -1896.1666

# Row 64 (BSM)
$ws.Range("H64").Value = 528.6667
$ws.Range("I64").Value = 558.5
$ws.Range("K64").Value = 558.5
$ws.Range("M64").Value = -333.5

# Row 67 (BSM)
$ws.Range("H67").Value = 528.6667
$ws.Range("I67").Value = 558.5
$ws.Range("K67").Value = 558.5
$ws.Range("M67").Value = 221.5

# Row 94 (BSM)
$ws.Range("H94").Value = 668.25
$ws.Range("I94").Value = 678
$ws.Range("J94").Value = 600
$ws.Range("K94").Value = 678
$ws.Range("L94").Value = 600
$ws.Range("M94").Value = -227
$ws.Range("N94").Value = -1502

# Row 105 (BSM)
$ws.Range("H105").Value = 2129.2222
$ws.Range("I105").Value = 1409.6923
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 1409.6923
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = 337.3077000000001
$ws.Range("N105").Value = -7494

# Row 107 (BSM)
$ws.Range("H107").Value = 69645.39999999999
$ws.Range("I107").Value = 3244.7693
$ws.Range("K107").Value = 3244.7693
$ws.Range("M107").Value = -1324.7693

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 13222.823
$ws.Range("I31").Value = 9978.799999999999
$ws.Range("K31").Value = 9978.799999999999
$ws.Range("M31").Value = -9683.799999999999

# Row 34 (CRP)
$ws.Range("H34").Value = 13222.823
$ws.Range("I34").Value = 9978.799999999999
$ws.Range("K34").Value = 9978.799999999999
$ws.Range("M34").Value = -9776.799999999999

# Row 50 (CRP)
$ws.Range("H50").Value = 65000
$ws.Range("J50").Value = 65000
$ws.Range("L50").Value = 65000
$ws.Range("N50").Value = -66250

# Row 51 (CRP)
$ws.Range("H51").Value = 44999.25
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 44999.25
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 44999.25
$ws.Range("N51").Value = -46471.25
$ws.Range("M51").ClearContents()

# Row 61 (CRP)
$ws.Range("H61").Value = 44999.25
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 44999.25
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 44999.25
$ws.Range("N61").Value = -45695.25
$ws.Range("M61").ClearContents()

# Row 132 (CRP)
$ws.Range("H132").Value = 1000000000
$ws.Range("I132").Value = 1000000000
$ws.Range("K132").Value = 3000000000
$ws.Range("M132").Value = -2999997470

$ws = $wb.Worksheets.Item("CUL")
# Row 2 (CUL)
$ws.Range("H2").Value = 632.5
$ws.Range("J2").Value = 697.7778
$ws.Range("L2").Value = 4186.6668
$ws.Range("N2").Value = -4412.6668

# Row 11 (CUL)
$ws.Range("H11").Value = 185077.94
$ws.Range("I11").Value = 198166.92
$ws.Range("K11").Value = 594500.76
$ws.Range("M11").Value = -594360.76

# Row 56 (CUL)
$ws.Range("H56").Value = 11806.794
$ws.Range("I56").Value = 11806.794
$ws.Range("K56").Value = 11806.794
$ws.Range("M56").Value = -11276.794

# Row 98 (CUL)
$ws.Range("H98").Value = 481.625
$ws.Range("I98").Value = 532.5
$ws.Range("J98").Value = 430.75
$ws.Range("K98").Value = 1597.5
$ws.Range("L98").Value = 1292.25
$ws.Range("M98").Value = -99.5
$ws.Range("N98").Value = -4288.25

# Row 105 (CUL)
$ws.Range("H105").Value = 10999
$ws.Range("J105").Value = 10999
$ws.Range("L105").Value = 32997
$ws.Range("N105").Value = -38239

# Row 114 (CUL)
$ws.Range("H114").Value = 201699.8
$ws.Range("J114").Value = 2666.6667
$ws.Range("L114").Value = 8000.000100000001
$ws.Range("N114").Value = -14508.0001

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (GSM)
$ws.Range("H70").Value = 11970.471
$ws.Range("I70").Value = 8920
$ws.Range("J70").Value = 16328.286
$ws.Range("K70").Value = 8920
$ws.Range("L70").Value = 16328.286
$ws.Range("M70").Value = -8650
$ws.Range("N70").Value = -16868.286

# Row 73 (GSM)
$ws.Range("H73").Value = 11970.471
$ws.Range("I73").Value = 8920
$ws.Range("J73").Value = 16328.286
$ws.Range("K73").Value = 8920
$ws.Range("L73").Value = 16328.286
$ws.Range("M73").Value = -7984
$ws.Range("N73").Value = -18200.286

# Row 107 (GSM)
$ws.Range("H107").Value = 1720.6666
$ws.Range("I107").Value = 1233.4706
$ws.Range("K107").Value = 1233.4706
$ws.Range("M107").Value = 686.5293999999999

# Row 123 (GSM)
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# Row 132 (GSM)
$ws.Range("H132").Value = 5003429
$ws.Range("I132").Value = 5437705.5
$ws.Range("J132").Value = 9249.5
$ws.Range("K132").Value = 16313116.5
$ws.Range("L132").Value = 27748.5
$ws.Range("M132").Value = -16310586.5
$ws.Range("N132").Value = -32808.5

$ws = $wb.Worksheets.Item("LTW")
# Row 57 (LTW)
$ws.Range("H57").Value = 26955
$ws.Range("J57").Value = 28950
$ws.Range("L57").Value = 28950
$ws.Range("N57").Value = -30082

# Row 132 (LTW)
$ws.Range("H132").Value = 53339520
$ws.Range("I132").Value = 60005708
$ws.Range("K132").Value = 180017124
$ws.Range("M132").Value = -180014594

$ws = $wb.Worksheets.Item("WVR")
# Row 81 (WVR)
$ws.Range("H81").Value = 168652.67
$ws.Range("I81").Value = 250479.25
$ws.Range("K81").Value = 500958.5
$ws.Range("M81").Value = -499897.5

# Row 84 (WVR)
$ws.Range("H84").Value = 168652.67
$ws.Range("I84").Value = 250479.25
$ws.Range("K84").Value = 2504792.5
$ws.Range("M84").Value = -2499488.5

# Row 122 (WVR)
$ws.Range("H122").Value = 2507.889
$ws.Range("I122").Value = 2040.8182
$ws.Range("K122").Value = 6122.4546
$ws.Range("M122").Value = -3672.4546

# Row 126 (WVR)
$ws.Range("H126").Value = 2045.4
$ws.Range("I126").Value = 2387
$ws.Range("J126").Value = 1703.8
$ws.Range("K126").Value = 7161
$ws.Range("L126").Value = 5111.4
$ws.Range("M126").Value = -4691
$ws.Range("N126").Value = -10051.4
